$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tokens = @(
    "Angel Warrior",
    "Bear",
    "Bird",
    "Cat",
    "Demon Berserker",
    "Dragon",
    "Dwarf Berserker",
    "Elf Warrior",
    "Foretell",
    "Giant Wizard",
    "Human Warrior",
    "Icy Manalith",
    "Kaya the Inexorable Emblem",
    "Koma's Coil",
    "Replicated Ring",
    "Shapeshifter",
    "Shard",
    "Spirit",
    "Tibalt, Cosmic Impostor Emblem",
    "Treasure",
    "Troll Warrior",
    "Tyvar Kell Emblem",
    "Zombie Berserker"
)

for ($i = 0; $i -lt $tokens.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $tokens[$i]
}
